$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("C1").Value = "Frecuencia del primer armonico"
$ws.Range("D1").Value = "Frecuencia del segundo armonico"
$ws.Range("E1").Value = "Frecuencia tercer armonico"

# Copy the header style (bold, centered, bordered) from C1 onto the new D1:E1 headers
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)

# Update data rows (new / revised frequency values for the three harmonics)
$ws.Cells.Item(2, 3).Value = 141.2991939307731
$ws.Cells.Item(2, 4).Value = 847.7951635846375
$ws.Cells.Item(2, 5).Value = 707.4442863916547
$ws.Cells.Item(3, 3).Value = 127.5949367088606
$ws.Cells.Item(3, 4).Value = 632.3037974683548
$ws.Cells.Item(3, 5).Value = 778.9367088607596
$ws.Cells.Item(4, 3).Value = 119.8338071287999
$ws.Cells.Item(4, 4).Value = 719.877542094905
$ws.Cells.Item(4, 5).Value = 239.6676142575989
$ws.Cells.Item(5, 3).Value = 130.5510968432318
$ws.Cells.Item(5, 4).Value = 257.2498662386306
$ws.Cells.Item(5, 5).Value = 781.1663991439273
$ws.Cells.Item(6, 3).Value = 772.9404769769444
$ws.Cells.Item(6, 4).Value = 128.1627799431863
$ws.Cells.Item(6, 5).Value = 889.2118649666377
$ws.Cells.Item(7, 3).Value = 575.763322586763
$ws.Cells.Item(7, 4).Value = 768.4430115683672
$ws.Cells.Item(7, 5).Value = 384.6007965105255
$ws.Cells.Item(8, 3).Value = 721.4938418752481
$ws.Cells.Item(8, 4).Value = 538.7365911799761
$ws.Cells.Item(8, 5).Value = 180.3734604688125
$ws.Cells.Item(9, 3).Value = 760.113710911874
$ws.Cells.Item(9, 4).Value = 607.0413295429698
$ws.Cells.Item(9, 5).Value = 912.3113929586707
$ws.Cells.Item(10, 3).Value = 962.8768226138291
$ws.Cells.Item(10, 4).Value = 1285.624832274801
$ws.Cells.Item(10, 5).Value = 1126.03989623401
$ws.Cells.Item(11, 3).Value = 285.2957746478878
$ws.Cells.Item(11, 4).Value = 568.3380281690143
$ws.Cells.Item(11, 5).Value = 856.788732394366
$ws.Cells.Item(12, 3).Value = 174.8206131767774
$ws.Cells.Item(12, 4).Value = 348.771472059143
$ws.Cells.Item(12, 5).Value = 693.1941726462278
$ws.Cells.Item(13, 3).Value = 638.122931916173
$ws.Cells.Item(13, 4).Value = 160.0320866339116
$ws.Cells.Item(13, 5).Value = 478.4919282061569
$ws.Cells.Item(14, 3).Value = 157.7424023154845
$ws.Cells.Item(14, 4).Value = 315.4848046309698
$ws.Cells.Item(14, 5).Value = 1420.40520984081
$ws.Cells.Item(15, 3).Value = 147.1032745591938
$ws.Cells.Item(15, 4).Value = 293.1989924433246
$ws.Cells.Item(15, 5).Value = 442.3173803526447
$ws.Cells.Item(16, 3).Value = 159.4117647058829
$ws.Cells.Item(16, 4).Value = 320
$ws.Cells.Item(16, 5).Value = 480
$ws.Cells.Item(17, 3).Value = 326.9007697985871
$ws.Cells.Item(17, 4).Value = 164.0830960666458
$ws.Cells.Item(17, 5).Value = 491.4056733101334
$ws.Cells.Item(18, 3).Value = 149.2472179794891
$ws.Cells.Item(18, 4).Value = 305.4767619463237
$ws.Cells.Item(18, 5).Value = 453.8511891773951
$ws.Cells.Item(19, 3).Value = 619.9973927779956
$ws.Cells.Item(19, 4).Value = 153.8260982922693
$ws.Cells.Item(19, 5).Value = 1239.99478555599
$ws.Cells.Item(20, 3).Value = 304.5260873060533
$ws.Cells.Item(20, 4).Value = 606.8011898062546
$ws.Cells.Item(20, 5).Value = 150.8159819921211
$ws.Cells.Item(21, 3).Value = 162.9268292682927
$ws.Cells.Item(21, 4).Value = 648.7804878048782
$ws.Cells.Item(21, 5).Value = 486.8292682926831
$ws.Cells.Item(22, 3).Value = 192.2869955156948
$ws.Cells.Item(22, 4).Value = 383.1390134529147
$ws.Cells.Item(22, 5).Value = 766.2780269058298
$ws.Cells.Item(23, 3).Value = 1376.100628930817
$ws.Cells.Item(23, 4).Value = 124.9475890985323
$ws.Cells.Item(23, 5).Value = 1490.146750524109
$ws.Cells.Item(24, 3).Value = 125.900153734121
$ws.Cells.Item(24, 4).Value = 252.7712598106641
$ws.Cells.Item(24, 5).Value = 1276.155028724007
$ws.Cells.Item(25, 3).Value = 1456.339535160198
$ws.Cells.Item(25, 4).Value = 1361.707186589788
$ws.Cells.Item(25, 5).Value = 1531.950306128515
$ws.Cells.Item(26, 3).Value = 1174.656469564136
$ws.Cells.Item(26, 4).Value = 131.5674936268106
$ws.Cells.Item(26, 5).Value = 393.7076416091522
$ws.Cells.Item(27, 3).Value = 166.0675772825307
$ws.Cells.Item(27, 4).Value = 661.3946800862686
$ws.Cells.Item(27, 5).Value = 332.1351545650614
$ws.Cells.Item(28, 3).Value = 1422.323635672181
$ws.Cells.Item(28, 4).Value = 850.3517779045451
$ws.Cells.Item(28, 5).Value = 994.8659440958359
$ws.Cells.Item(29, 3).Value = 1316.692667706708
$ws.Cells.Item(29, 4).Value = 1170.046801872075
$ws.Cells.Item(29, 5).Value = 732.6053042121684
$ws.Cells.Item(30, 3).Value = 1427.76432606941
$ws.Cells.Item(30, 4).Value = 1319.612590799032
$ws.Cells.Item(30, 5).Value = 1049.233252623083
$ws.Cells.Item(31, 3).Value = 316.792265036258
$ws.Cells.Item(31, 4).Value = 159.8180008531217
$ws.Cells.Item(31, 5).Value = 791.1275415896489
$ws.Cells.Item(32, 3).Value = 227.1823204419888
$ws.Cells.Item(32, 4).Value = 449.9447513812156
$ws.Cells.Item(32, 5).Value = 678.0110497237565
$ws.Cells.Item(33, 3).Value = 246.8429830831546
$ws.Cells.Item(33, 4).Value = 1729.807005003574
$ws.Cells.Item(33, 5).Value = 494.639027877055
$ws.Cells.Item(34, 3).Value = 1552.255326999354
$ws.Cells.Item(34, 4).Value = 1032.746056636841
$ws.Cells.Item(34, 5).Value = 1293.23863112259
$ws.Cells.Item(35, 3).Value = 248.4276729559747
$ws.Cells.Item(35, 4).Value = 494.7589098532499
$ws.Cells.Item(35, 5).Value = 743.1865828092241
$ws.Cells.Item(36, 3).Value = 971.0910786447002
$ws.Cells.Item(36, 4).Value = 728.6291576002486
$ws.Cells.Item(36, 5).Value = 1458.501709667393
$ws.Cells.Item(37, 3).Value = 463.77738685431
$ws.Cells.Item(37, 4).Value = 694.0668479130018
$ws.Cells.Item(37, 5).Value = 233.4879257956181
$ws.Cells.Item(38, 3).Value = 668.3221000354733
$ws.Cells.Item(38, 4).Value = 222.7740333451579
$ws.Cells.Item(38, 5).Value = 445.5480666903159
$ws.Cells.Item(39, 3).Value = 708.2018927444788
$ws.Cells.Item(39, 4).Value = 176.1303890641429
$ws.Cells.Item(39, 5).Value = 1240.273396424815
$ws.Cells.Item(40, 3).Value = 181.182595760506
$ws.Cells.Item(40, 4).Value = 359.6876162142062
$ws.Cells.Item(40, 5).Value = 541.167720342135
$ws.Cells.Item(41, 3).Value = 172.9544806614081
$ws.Cells.Item(41, 4).Value = 346.289081060534
$ws.Cells.Item(41, 5).Value = 520.3839209350945
$ws.Cells.Item(42, 3).Value = 782.8171445009111
$ws.Cells.Item(42, 4).Value = 263.4960207114773
$ws.Cells.Item(42, 5).Value = 522.7730367245176
$ws.Cells.Item(43, 3).Value = 273.0348258706472
$ws.Cells.Item(43, 4).Value = 1095.323383084577
$ws.Cells.Item(43, 5).Value = 548.4577114427866
$ws.Cells.Item(44, 3).Value = 907.1871505540307
$ws.Cells.Item(44, 4).Value = 226.4918166107554
$ws.Cells.Item(44, 5).Value = 681.1019619802782
$ws.Cells.Item(45, 3).Value = 177.3352299668086
$ws.Cells.Item(45, 4).Value = 350.8771929824561
$ws.Cells.Item(45, 5).Value = 700.806069227122
$ws.Cells.Item(46, 3).Value = 293.4455802766088
$ws.Cells.Item(46, 4).Value = 584.4858689116054
$ws.Cells.Item(46, 5).Value = 877.9314491882142
$ws.Cells.Item(47, 3).Value = 1069.073783359498
$ws.Cells.Item(47, 4).Value = 802.9827315541602
$ws.Cells.Item(47, 5).Value = 534.5368916797488
$ws.Cells.Item(48, 3).Value = 280.620842572062
$ws.Cells.Item(48, 4).Value = 841.1529933481161
$ws.Cells.Item(48, 5).Value = 562.3059866962312
$ws.Cells.Item(49, 3).Value = 842.0334356874791
$ws.Cells.Item(49, 4).Value = 1122.25634027067
$ws.Cells.Item(49, 5).Value = 560.9007164790173
$ws.Cells.Item(50, 3).Value = 801.0645375914837
$ws.Cells.Item(50, 4).Value = 199.6007984031939
$ws.Cells.Item(50, 5).Value = 398.5362608117102
$ws.Cells.Item(51, 3).Value = 187.5945537065054
$ws.Cells.Item(51, 4).Value = 375.8614893259373
$ws.Cells.Item(51, 5).Value = 747.6886871743145
$ws.Cells.Item(52, 3).Value = 198.4196663740122
$ws.Cells.Item(52, 4).Value = 398.0099502487565
$ws.Cells.Item(52, 5).Value = 801.287679250805
$ws.Cells.Item(53, 3).Value = 813.8876658909367
$ws.Cells.Item(53, 4).Value = 203.534813510284
$ws.Cells.Item(53, 5).Value = 1023.209006855777
$ws.Cells.Item(54, 3).Value = 205.5563552612639
$ws.Cells.Item(54, 4).Value = 822.8012091550308
$ws.Cells.Item(54, 5).Value = 1028.357564416295
